# Adds the EXP2/EXP1 ratio column (J) and a new scatter chart plotting it
# against M (column A), mirroring the author's commit:
#   "Added exponential ratio, updated xlsx file"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column: header + ratio formulas (EXP SECOND / EXP FIRST) -----
$ws.Range("J20").Value = "EXP2/EXP1"
$ws.Range("J21").Formula = "=(H21/D21)"
$ws.Range("J22:J34").Formula = "=(H22/D22)"

# --- 2. New scatter chart (EXP2/EXP1 vs M, log-scale X axis) -------------
$co = $ws.ChartObjects().Add(695.24, 10.92, 321.44, 213.19)
$co.Name = "Chart 1"

$chart = $co.Chart
$chart.ChartType = 74          # xlXYScatterLines

# Configure the category (X) axis as log-base-10 *before* wiring up the
# series data -- the axis scaling only "sticks" when it is the chart's very
# first axis touch.
$axX = $chart.Axes(1)
$axX.LogBase = 10
$axX.HasMajorGridlines = $true

$axY = $chart.Axes(2)
$axY.HasMajorGridlines = $true

$chart.SeriesCollection().NewSeries()
$ser = $chart.SeriesCollection(1)
$ser.Name = "=Sheet1!`$J`$20"
$ser.XValues = $ws.Range("A21:A34")
$ser.Values = $ws.Range("J21:J34")
$ser.MarkerStyle = 8           # xlMarkerStyleCircle
$ser.MarkerSize = 5

# --- 3. Hidden "_xlchart" defined names Excel stamps per inserted chart --
$wb.Names.Add("_xlchart.v1.0", "=Sheet1!`$A`$20").Visible = $false
$wb.Names.Add("_xlchart.v1.1", "=Sheet1!`$A`$21:`$A`$34").Visible = $false
$wb.Names.Add("_xlchart.v1.2", "=Sheet1!`$J`$20").Visible = $false
$wb.Names.Add("_xlchart.v1.3", "=Sheet1!`$J`$21:`$J`$34").Visible = $false
$wb.Names.Add("_xlchart.v2.4", "=Sheet1!`$A`$20").Visible = $false
$wb.Names.Add("_xlchart.v2.5", "=Sheet1!`$A`$21:`$A`$34").Visible = $false
$wb.Names.Add("_xlchart.v2.6", "=Sheet1!`$J`$20").Visible = $false
$wb.Names.Add("_xlchart.v2.7", "=Sheet1!`$J`$21:`$J`$34").Visible = $false

# --- 4. Selection moved by the author while working near the new chart ---
$ws.Range("P19").Select()
